# Generate Report for Handoff
#
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the associated timestamps on the Overview,
# zh-cn and de-de sheets of the handback report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-04 06:42:48"

# --- zh-cn sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-04 06:42:44"

# --- de-de sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-04 06:42:48"

# --- Column widths widen to fit "Ready for handoff" -------------------
# Overview: columns E (zh-cn) and F (de-de) status columns.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336

# zh-cn / de-de: column C is the Status column.
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333336
$dede.Columns.Item(3).ColumnWidth = 16.333333333333336
